# Final push for last endpoint documentation
# Inserts a new "GET /v1/orders/amazon/creds" endpoint row into the
# "/v1/orders/amazon" block on Sheet1 (between the existing GET row and
# the POST row), pushing the POST/Delete/PUT rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 28 (shifts old rows 28-30 down to 29-31)
$ws.Rows("28:28").Insert()

# New row 28: GET /v1/orders/amazon/creds
$ws.Range("C28").Value = "/v1/orders/amazon/creds"
$ws.Range("D28").Value = "GET"
$ws.Range("E28").Value = "store_login.checkForAmazonCredentials_v1"

# The path column for the rest of this endpoint group now reads
# "/v1/orders/amazon/creds" as well (rows 29-31, formerly 28-30)
$ws.Range("C29").Value = "/v1/orders/amazon/creds"
$ws.Range("C30").Value = "/v1/orders/amazon/creds"
$ws.Range("C31").Value = "/v1/orders/amazon/creds"

# Restore the selection to match the author's final cursor position
$ws.Range("G23").Select() | Out-Null
